$d = $word.ActiveDocument

$d.Content.Find.Execute("28 de junio de 2017.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "06 de julio de 2017.", 2)

$d.Content.Find.Execute("Fiscalia Nacional en lo Criminal y Correccional Nro 1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fiscalia Nacional en lo Criminal y Correccional Nro 3", 2)

$d.Content.Find.Execute('$ 942,00', $true, $false, $false, $false, $false,
                         $true, 1, $false, '$ 200,00', 2)

$d.Content.Find.Execute("Finalmente, la presente erogación de fondos es solicitada por este curso debido a que Es un pedido urgente", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Finalmente, la presente erogación de fondos es solicitada por este curso debido a que asdfasd", 2)
